$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value2 = 'Administrator, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Amira Sobhy'
$ws.Range('G3').Value2 = 'Dr. Asmaa Reda, Administrator, Dr. Majorelle Magdy, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud'
$ws.Range('G4').Value2 = 'Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud'
$ws.Range('G5').Value2 = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi'
$ws.Range('G6').Value2 = 'Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany'
$ws.Range('G7').Value2 = 'Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Amera Ahmad Saad, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range('G8').Value2 = 'Dr. Nada Mohammad, Dr. Abeer Ragab'
$ws.Range('G9').Value2 = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G11').Value2 = 'Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad'
$ws.Range('G12').Value2 = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Dina Adel'
$ws.Range('G13').Value2 = 'Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh'
$ws.Range('G15').Value2 = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range('G19').Value2 = 'Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range('G20').Value2 = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range('G25').Value2 = 'Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud'
$ws.Range('G27').Value2 = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range('G28').Value2 = 'Dr. Aya Emad, Dr. Maryam Ashraf'
$ws.Range('G30').Value2 = 'Dr. Shorok Mohammad, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy'
